$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style, borders, fill) of column L (reference_id column)
# over to the new column M so the new "transaction status" / "link" columns
# look consistent with the rest of the table (header style, data-row style,
# and the styled-but-empty placeholder rows 4-10).
$ws.Range("L1:L10").Copy()
$ws.Range("M1:M10").PasteSpecial(-4122)

# New header + data cells for the additional columns.
$ws.Range("M1").Value = "transaction_status_blockchain"
$ws.Range("M2").Value = "transaction_status_blockchain_link"
$ws.Range("M3").Value = "www_link"

# Widen the new column so the long header/values are readable.
$ws.Columns.Item(13).ColumnWidth = 41.5
